$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: "Roller er sterkere men ball ..." - merge the three runs
# (which were split apart only so Word could flag "sterkere" with a
# gramStart/gramEnd proofing mark) back into a single run of text.
# A contained Find/Replace merges the runs and keeps the surrounding
# run's formatting (confirmed experimentally against this runtime).
# -----------------------------------------------------------------
$d.Content.Find.Execute(
    "Roller er sterkere men ball", $true, $false, $false, $false, $false,
    $true, 1, $false, "Roller er sterkere men ball", 2) | Out-Null

# -----------------------------------------------------------------
# Change 2: add new content at the end of the document.
#   Paragraph A (the existing empty last paragraph) gets the
#       "Posisjonering av kablene ..." text.
#   Paragraph B: a new empty paragraph.
#   Paragraph C: a new Heading1 paragraph "Albue".
#   Paragraph D: a new paragraph with the "Hvor langt kan armen ..."
#       text.
# -----------------------------------------------------------------

# --- Paragraph A: fill the existing trailing empty paragraph -----
$parA = $d.Paragraphs.Last
$parA.Range.InsertAfter(
    "Posisjonering av kablene er utrolig vanskelig. Det må tenkes på mens man designer armen, hvis ikke blir designet umulig uten sleperinger. Egentlig burde signalkablene ligge separat fra supply spenning, men det er helt umulig å få til uten å legge noe kabler synlig.")
$parA.Range.LanguageID = "nb-NO"

# --- split off a new empty Paragraph B after A --------------------
$d.Content.Find.Execute(
    "legge noe kabler synlig.", $true, $false, $false, $false, $false,
    $true, 1, $false, "legge noe kabler synlig.^p", 2) | Out-Null

# --- Paragraph C: new Heading1 paragraph "Albue" ------------------
$parB = $d.Paragraphs.Last
$parB.Range.InsertParagraphAfter()

$parC = $d.Paragraphs.Last
$parC.Range.InsertAfter("Albue")
$parC.Range.LanguageID = "nb-NO"

# split off a new empty Paragraph D after C, while C is still Normal
# styled, then turn C into Heading1 - this avoids the new paragraph
# (D) inheriting the Heading1 style.
$d.Content.Find.Execute(
    "Albue", $true, $false, $false, $false, $false,
    $true, 1, $false, "Albue^p", 2) | Out-Null

$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $cand = $paras.Item($i)
    if ($cand.Range.Text.Trim() -eq "Albue") {
        $cand.Style = "Heading 1"
    }
}

# --- Paragraph D: fill the new trailing empty paragraph -----------
$parD = $d.Paragraphs.Last
$parD.Range.InsertAfter(
    "Hvor langt kan armen legges sammen i albuen? Dette er helt avhengig av geometrien på armen og plassering av worm-gir motor. Vi vil gjerne at armen skal se rett ut i utstrakt posisjon, derfor ønsker vi ")
$parD.Range.LanguageID = "nb-NO"

# -----------------------------------------------------------------
# Now that all text/paragraph-structure edits are finished (further
# Find/Replace or paragraph splits would otherwise re-merge runs
# with identical formatting), split "supply" and "worm" off into
# their own runs - this mirrors the run boundaries the proofErr
# (spell-check) marks delimit in the source diff. Toggling Bold on
# and back off forces the runtime to break the run without leaving
# any residual character formatting behind.
# -----------------------------------------------------------------
$parAFinal = $null
$parDFinal = $null
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $cand = $paras.Item($i)
    if ($cand.Range.Text.Contains("signalkablene ligge separat fra supply")) {
        $parAFinal = $cand
    }
    if ($cand.Range.Text.Contains("plassering av worm-gir motor")) {
        $parDFinal = $cand
    }
}

$fullA = $parAFinal.Range.Text
$supplyStart = $parAFinal.Range.Start + $fullA.IndexOf("supply")
$supplyEnd = $supplyStart + "supply".Length
$supplyRange = $d.Range($supplyStart, $supplyEnd)
$supplyRange.Bold = $true
$supplyRange.Bold = $false

$fullD = $parDFinal.Range.Text
$wormStart = $parDFinal.Range.Start + $fullD.IndexOf("worm")
$wormEnd = $wormStart + "worm".Length
$wormRange = $d.Range($wormStart, $wormEnd)
$wormRange.Bold = $true
$wormRange.Bold = $false

Write-Host "Done."
